# Applies the diff:
#  - For data rows 2..99 (except row 36, which is unchanged), decrement
#    column E ("剩余") by 1.
#  - Row 95 is a special case: E95 becomes 10 (instead of 1-1=0) and
#    F95 ("开始时间") changes from 20260120 to 20260130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    $cell = $ws.Cells.Item($row, 5)  # Column E
    if ($row -eq 95) {
        $cell.Value2 = 10
        $ws.Cells.Item($row, 6).Value2 = 20260130  # Column F
    } else {
        $current = $cell.Value2
        $cell.Value2 = $current - 1
    }
}
